# Auto-generated Excel COM-interop script applying numeric corrections
# to the Leve profit-tracking tables (H/I/J/K/L/M/N columns) across all 8
# crafting-job sheets, per the scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 274.5
$ws.Range("I6").Value = 250
$ws.Range("K6").Value = 750
$ws.Range("M6").Value = -638
$ws.Range("H19").Value = 3225.6
$ws.Range("J19").Value = 3702.3076
$ws.Range("L19").Value = 3702.3076
$ws.Range("N19").Value = -4052.3076
$ws.Range("H62").Value = 4548.2383
$ws.Range("I62").Value = 4212.7896
$ws.Range("J62").Value = 7735
$ws.Range("K62").Value = 4212.7896
$ws.Range("L62").Value = 7735
$ws.Range("M62").Value = -3588.7896
$ws.Range("N62").Value = -8983
$ws.Range("H65").Value = 4548.2383
$ws.Range("I65").Value = 4212.7896
$ws.Range("J65").Value = 7735
$ws.Range("K65").Value = 21063.948
$ws.Range("L65").Value = 38675
$ws.Range("M65").Value = -17943.948
$ws.Range("N65").Value = -44915
$ws.Range("H74").Value = 7099.273
$ws.Range("I74").Value = 7099.273
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 7099.273
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -6163.273
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 7099.273
$ws.Range("I77").Value = 7099.273
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 35496.365
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -30816.365
$ws.Range("N77").ClearContents()
$ws.Range("H100").Value = 2365.5715
$ws.Range("I100").Value = 826
$ws.Range("J100").Value = 2981.4
$ws.Range("K100").Value = 826
$ws.Range("L100").Value = 2981.4
$ws.Range("M100").Value = -285
$ws.Range("N100").Value = -4063.4
$ws.Range("H137").Value = 8446.426
$ws.Range("I137").Value = 10620.883
$ws.Range("K137").Value = 31862.649
$ws.Range("M137").Value = -29312.649

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 38681.68
$ws.Range("I32").Value = 41311.04
$ws.Range("J32").Value = 4500
$ws.Range("K32").Value = 41311.04
$ws.Range("L32").Value = 4500
$ws.Range("M32").Value = -41024.04
$ws.Range("N32").Value = -5074
$ws.Range("H43").Value = 59159.5
$ws.Range("J43").Value = 92377
$ws.Range("L43").Value = 92377
$ws.Range("N43").Value = -93003
$ws.Range("H45").Value = 2487.6924
$ws.Range("I45").Value = 1771.3889
$ws.Range("K45").Value = 1771.3889
$ws.Range("M45").Value = -1394.3889
$ws.Range("H122").Value = 2479.4814
$ws.Range("I122").Value = 2400.0476
$ws.Range("K122").Value = 7200.1428
$ws.Range("M122").Value = -4750.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 59.95238
$ws.Range("I7").Value = 71.71429
$ws.Range("K7").Value = 71.71429
$ws.Range("M7").Value = 41.28570999999999
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 7692585
$ws.Range("I7").Value = 299.66666
$ws.Range("J7").Value = 19231012
$ws.Range("K7").Value = 898.9999799999999
$ws.Range("L7").Value = 57693036
$ws.Range("M7").Value = -786.9999799999999
$ws.Range("N7").Value = -57693260
$ws.Range("H92").Value = 1828.5834
$ws.Range("I92").Value = 1794.3
$ws.Range("J92").Value = 2000
$ws.Range("K92").Value = 5382.9
$ws.Range("L92").Value = 6000
$ws.Range("M92").Value = -4134.9
$ws.Range("N92").Value = -8496
$ws.Range("H107").Value = 1894.3182
$ws.Range("I107").Value = 2927.5454
$ws.Range("J107").Value = 861.0909
$ws.Range("K107").Value = 8782.6362
$ws.Range("L107").Value = 2583.2727
$ws.Range("M107").Value = -6862.636200000001
$ws.Range("N107").Value = -6423.2727
$ws.Range("H131").Value = 1467.0834
$ws.Range("I131").Value = 1057.2667
$ws.Range("J131").Value = 2150.111
$ws.Range("K131").Value = 3171.800099999999
$ws.Range("L131").Value = 6450.333
$ws.Range("M131").Value = 1868.199900000001
$ws.Range("N131").Value = -16530.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 23625.322
$ws.Range("I102").Value = 30705.762
$ws.Range("J102").Value = 2384
$ws.Range("K102").Value = 30705.762
$ws.Range("L102").Value = 2384
$ws.Range("M102").Value = -29083.762
$ws.Range("N102").Value = -5628
$ws.Range("H105").Value = 330000
$ws.Range("J105").Value = 600000
$ws.Range("L105").Value = 600000
$ws.Range("N105").Value = -606988

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 894.4167
$ws.Range("I16").Value = 866.5455
$ws.Range("J16").Value = 1201
$ws.Range("K16").Value = 866.5455
$ws.Range("L16").Value = 1201
$ws.Range("M16").Value = -696.5455
$ws.Range("N16").Value = -1541
$ws.Range("H61").Value = 2450
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H113").Value = 2450
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 3827.7917
$ws.Range("I122").Value = 3163.8125
$ws.Range("J122").Value = 5155.75
$ws.Range("K122").Value = 9491.4375
$ws.Range("L122").Value = 15467.25
$ws.Range("M122").Value = -7041.4375
$ws.Range("N122").Value = -20367.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10280.4375
$ws.Range("I81").Value = 11998.917
$ws.Range("J81").Value = 5125
$ws.Range("K81").Value = 23997.834
$ws.Range("L81").Value = 10250
$ws.Range("M81").Value = -22936.834
$ws.Range("N81").Value = -12372
$ws.Range("H84").Value = 10280.4375
$ws.Range("I84").Value = 11998.917
$ws.Range("J84").Value = 5125
$ws.Range("K84").Value = 119989.17
$ws.Range("L84").Value = 51250
$ws.Range("M84").Value = -114685.17
$ws.Range("N84").Value = -61858
$ws.Range("H132").Value = 2511.8333
$ws.Range("I132").Value = 1559.75
$ws.Range("K132").Value = 4679.25
$ws.Range("M132").Value = -2149.25
